$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 112243463
$ws.Range("B2").Value = 77389
$ws.Range("E2").Value = 228912
$ws.Range("F2").Value = "Mörk kolflarnlav"
$ws.Range("G2").Value = "Carbonicola myrmecina"
$ws.Range("H2").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q2").Value = 410608
$ws.Range("R2").Value = 6710914

# Row 3
$ws.Range("A3").Value = 112243462
$ws.Range("B3").Value = 78228
$ws.Range("E3").Value = 6453
$ws.Range("F3").Value = "Vedskivlav"
$ws.Range("G3").Value = "Hertelidea botryosa"
$ws.Range("H3").Value = "(Fr.) Printzen & Kantvilas"

# Row 4
$ws.Range("A4").Value = 112243469
$ws.Range("B4").Value = 77636
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("Q4").Value = 410486
$ws.Range("R4").Value = 6710828

# Row 5
$ws.Range("B5").Value = 77389

# Row 6
$ws.Range("A6").Value = 112243468
$ws.Range("B6").Value = 77636
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = "Garnlav"
$ws.Range("G6").Value = "Alectoria sarmentosa"
$ws.Range("H6").Value = "(Ach.) Ach."
$ws.Range("Q6").Value = 410566
$ws.Range("R6").Value = 6710872

# Row 7
$ws.Range("A7").Value = 112243461
$ws.Range("B7").Value = 77389
$ws.Range("E7").Value = 228912
$ws.Range("F7").Value = "Mörk kolflarnlav"
$ws.Range("G7").Value = "Carbonicola myrmecina"
$ws.Range("H7").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q7").Value = 410598
$ws.Range("R7").Value = 6710899

# Row 8
$ws.Range("A8").Value = 112243460
$ws.Range("B8").Value = 78228
$ws.Range("E8").Value = 6453
$ws.Range("F8").Value = "Vedskivlav"
$ws.Range("G8").Value = "Hertelidea botryosa"
$ws.Range("H8").Value = "(Fr.) Printzen & Kantvilas"
